{"js": "const body = context.document.body;\n\n// The source edit rewrites a single sentence:\n//   \"...o site contar\u00e1 com resenhas de livros, recomenda\u00e7\u00f5es de leitura e...\"\n// becomes\n//   \"...o site contar\u00e1 com a import\u00e2ncia dos livros, benef\u00edcios de leitura e...\"\n//\n// Search on the full unique phrase (not just the changed words) so this\n// can't accidentally match the similar \"recomenda\u00e7\u00f5es de leitura\" /\n// \"recomenda\u00e7\u00f5es\" phrasing that appears elsewhere in the document.\nconst target =\n  \"contar\u00e1 com resenhas de livros, recomenda\u00e7\u00f5es de leitura e dicas\";\nconst replacement =\n  \"contar\u00e1 com a import\u00e2ncia dos livros, benef\u00edcios de leitura e dicas\";\n\nconst results = body.search(target, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# The source edit rewrites a single sentence:\n#   \"...o site contar\u00e1 com resenhas de livros, recomenda\u00e7\u00f5es de leitura e...\"\n# becomes\n#   \"...o site contar\u00e1 com a import\u00e2ncia dos livros, benef\u00edcios de leitura e...\"\n#\n# Search on the full unique phrase (not just the changed words) so this\n# can't accidentally match the similar \"recomenda\u00e7\u00f5es de leitura\" /\n# \"recomenda\u00e7\u00f5es\" phrasing that appears elsewhere in the document, and\n# replace only the first (only) occurrence.\n$find = \"contar\u00e1 com resenhas de livros, recomenda\u00e7\u00f5es de leitura e dicas\"\n$replace = \"contar\u00e1 com a import\u00e2ncia dos livros, benef\u00edcios de leitura e dicas\"\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$r.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 1)\n"}
